$d = $word.ActiveDocument

# The header contains a REF field to bookmark "Schulhalbjahr12neu" that
# displays the school year (e.g. "2016/2017"). Replace the whole field
# (begin/instrText/separate/result/end) with a literal "${schuljahr}"
# placeholder, matching the other "${...}" template placeholders already
# used throughout the document (${name}, ${klasse}, ...).

$hdr = $d.Sections.Item(1).Headers.Item(1)
$hdr.Range.Find.Execute("2016/2017", $false, $false, $false, $false, $false,
                         $true, 1, $false, '${schuljahr}', 2) | Out-Null
